$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage for numeric-looking strings,
# so Excel does not silently convert them to numbers (losing formatting like trailing zeros,
# or multi-dot grouped numbers), while keeping the cell unstyled (same as original).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '66.929.37'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '3.825.40'
$ws.Range('E3').Value = '  +3.52%  '
$ws.Range('E4').Value = '  +0.07%  '
Set-TextValue 'D5' '412.67'
$ws.Range('E5').Value = '  -1.62%  '
Set-TextValue 'D6' '132.41'
$ws.Range('E6').Value = '  +1.50%  '
$ws.Range('D7').Value = '3.820.79'
$ws.Range('E7').Value = '  +3.59%  '
Set-TextValue 'D8' '0.615'
$ws.Range('E8').Value = '  -4.25%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  -2.79%  '
Set-TextValue 'D11' '0.170'
$ws.Range('E11').Value = '  -5.58%  '
Set-TextValue 'D12' '0.0000376'
$ws.Range('E12').Value = '  -3.62%  '
Set-TextValue 'D13' '41.11'
$ws.Range('E13').Value = '  -5.09%  '
$ws.Range('D14').Value = '4.432.65'
$ws.Range('E14').Value = '  +3.55%  '
Set-TextValue 'D15' '10.03'
$ws.Range('E15').Value = '  -5.81%  '
Set-TextValue 'D16' '14.91'
$ws.Range('E16').Value = '  +13.63%  '
$ws.Range('E17').Value = '  -1.20%  '
$ws.Range('D18').Value = '3.831.42'
$ws.Range('E18').Value = '  +0.55%  '
Set-TextValue 'D19' '19.52'
$ws.Range('E19').Value = '  -4.83%  '
$ws.Range('D20').Value = '67.260.14'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('E21').Value = '  -3.29%  '
Set-TextValue 'D22' '416.16'
$ws.Range('E22').Value = '  -6.12%  '
Set-TextValue 'D23' '14.82'
$ws.Range('E23').Value = '  -8.10%  '
Set-TextValue 'D24' '86.13'
$ws.Range('E24').Value = '  -4.25%  '
$ws.Range('E25').Value = '  -1.80%  '
Set-TextValue 'D26' '36.82'
$ws.Range('E26').Value = '  -1.61%  '
$ws.Range('E27').Value = '  +13.69%  '
$ws.Range('E28').Value = '  -4.58%  '
$ws.Range('E29').Value = '  -6.92%  '
Set-TextValue 'D30' '699.32'
$ws.Range('E30').Value = '  +7.07%  '
Set-TextValue 'D31' '0.123'
$ws.Range('E31').Value = '  -1.60%  '
$ws.Range('E32').Value = '  -2.04%  '
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('E34').Value = '  -0.69%  '
Set-TextValue 'D35' '0.153'
$ws.Range('E35').Value = '  -8.07%  '
Set-TextValue 'D36' '39.14'
$ws.Range('E36').Value = '  -6.75%  '
Set-TextValue 'D37' '1.00'
$ws.Range('E37').Value = '  +0.07%  '
Set-TextValue 'D38' '55.55'
$ws.Range('E38').Value = '  -2.84%  '
$ws.Range('D39').Value = '0.0₃0777'
$ws.Range('E39').Value = '  +6.17%  '
Set-TextValue 'D42' '0.997'
$ws.Range('E42').Value = '  -0.11%  '
Set-TextValue 'D43' '27.72'
$ws.Range('E43').Value = '  -5.20%  '
$ws.Range('E44').Value = '  -8.80%  '
Set-TextValue 'D45' '148.82'
$ws.Range('E45').Value = '  -0.04%  '
Set-TextValue 'D46' '3.19'
$ws.Range('E46').Value = '  +19.15%  '
Set-TextValue 'D47' '3.34'
$ws.Range('E47').Value = '  -2.31%  '
Set-TextValue 'D48' '4.43'
$ws.Range('E48').Value = '  +1.98%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('E50').Value = '  -1.29%  '
$ws.Range('E51').Value = '  -2.04%  '

# Row 40/41 coin swap: ThetaToken moves to row 40, VeChain moves to row 41
$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D40' '3.09'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D41' '0.0462'
$ws.Range('E41').Value = '  -6.42%  '
